$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.065.86'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.37%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.417.40'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.11%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.10%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  +1.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.416.66'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.15%  '
$ws.Range('E10').Value = '  +4.88%  '
$ws.Range('E11').Value = '  +0.81%  '
$ws.Range('E12').Value = '  +2.12%  '
$ws.Range('E13').Value = '  +4.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.25'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.25%  '
$ws.Range('E15').Value = '  +9.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.853.31'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.882.61'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.415.17'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.09'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '324.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.66%  '
$ws.Range('E22').Value = '  +3.08%  '
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.77'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.82'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.85%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.23'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '563.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +13.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.541.42'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.35'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0935'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +9.40%  '
$ws.Range('E32').Value = '  +6.06%  '
$ws.Range('E33').Value = '  +2.79%  '
$ws.Range('E34').Value = '  +4.40%  '
$ws.Range('E35').Value = '  +3.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.74'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.71%  '
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.97'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.96%  '
$ws.Range('E39').Value = '  +5.55%  '
$ws.Range('E40').Value = '  +2.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.81'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '147.07'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.20%  '
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.47'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '151.97'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.44%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +12.00%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.63'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.78%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0546'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.24%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.41'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.59%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.590'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.36%  '
$ws.Range('E51').Value = '  +3.79%  '
